# Refreshes the cryptos list with updated Price/Volume(1h) figures pulled
# from the latest coinranking.com snapshot, and re-syncs row order for the
# handful of coins whose new figures swapped their relative rank
# (Chainlink <-> WrappedEther, Stellar <-> ThetaToken, NEARProtocol <-> WEMIXToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.734.96"
$ws.Range("E2").Value = "  +3.52%  "

$ws.Range("D3").Value = "3.760.89"
$ws.Range("E3").Value = "  +7.19%  "

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "421.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.22"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "3.746.03"
$ws.Range("E7").Value = "  +6.99%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.651"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.72%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.775"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.187"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000429"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +60.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.60%  "

$ws.Range("D15").Value = "4.362.34"
$ws.Range("E15").Value = "  +7.39%  "

$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.756.03"
$ws.Range("E17").Value = "  +6.70%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.97%  "

$ws.Range("E20").Value = "  +2.82%  "

$ws.Range("D21").Value = "67.723.10"
$ws.Range("E21").Value = "  +3.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "451.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +19.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "90.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.67%  "

$ws.Range("E29").Value = "  +5.65%  "

$ws.Range("E30").Value = "  +6.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.76%  "

$ws.Range("E32").Value = "  +0.72%  "

$ws.Range("E33").Value = "  -3.20%  "

$ws.Range("E34").Value = "  +1.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "42.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("E38").Value = "  -2.51%  "

$ws.Range("D39").Value = "0.0₃0752"
$ws.Range("E39").Value = "  +1.26%  "

$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +30.45%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.149"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "27.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +27.81%  "

$ws.Range("E44").Value = "  +3.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.47%  "

$ws.Range("E48").Value = "  +22.68%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.93%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.43%  "

$ws.Range("E51").Value = "  -1.52%  "
